# Convert the year-header row (E1:BL1) from text labels like "1960 [YR1960]"
# into plain numeric years (1960-2019), left-aligned, so the years can be used
# numerically (e.g. for correlation / loop code over the columns).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headerRange = $ws.Range("E1:BL1")

# Apply left alignment to the header range first; Excel then carries that
# style (xfId) onto each cell as its literal value is written below.
$headerRange.HorizontalAlignment = -4131   # xlLeft

$startYear = 1960
for ($i = 0; $i -lt 60; $i++) {
    $col = 5 + $i   # column E = 5
    $ws.Cells.Item(1, $col).Value = $startYear + $i
}

# Update the selection to mirror the edited range.
$headerRange.Select() | Out-Null
